# Updated non-tracing input data
#
# The total population inputs for the "ideal-format" sheet (B3 = Germany,
# C3 = total N used in the model) used to be hard-coded numbers. They are
# replaced with formulas that sum the traced-through sub-population rows
# (S, E, I_asym, I_sym, I_sev, R, D -> rows 16-22), so the totals stay in
# sync with the underlying compartment numbers instead of drifting out of
# date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"
$ws.Range("C3").Formula = "=C16+C17+C18+C19+C20+C21+C22"

# Selection moved from A11:XFD11 (whole row) to B1 as the active cell.
$ws.Range("B1").Select()
